$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.779.55"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.340.88"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -9.08%  "
$ws.Range("D9").Value = "2.339.99"
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "2.764.77"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "60.762.12"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "2.337.86"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -6.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.07%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "2.456.28"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "495.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("E32").Value = "  -6.85%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "142.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.35%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  -6.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  -6.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.566"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("E51").Value = "  -2.85%  "
